# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.574.46"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "2.430.54"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.27"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").Value = "2.441.18"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("D14").Value = "2.863.74"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").Value = "57.494.68"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.75"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").Value = "2.440.75"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.42"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "315.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.09"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.63"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.43"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.405"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.158"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.21"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.26"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0721"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.22"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.67"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.91"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.26"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.45"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.774"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "272.13"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.37"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.88"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.583"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0907"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.34"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0482"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.09"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0210"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.49"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.51%  "
